$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ===========================================================================
# 1) Re-synchronise betting-odds rows whose fixtures were re-ordered by the
#    upstream scraper re-run on 05-11-2023 (same underlying matches, moved to
#    a different row). Only columns F:V (home..url) change; columns A:E
#    (Indice/pais/torneio/temporada/data_partida) are untouched.
# ===========================================================================
$ws.Range("F48").Value = 'Molde'
$ws.Range("G48").Value = 3
$ws.Range("H48").Value = 'Stromsgodset'
$ws.Range("I48").Value = 2
$ws.Range("J48").Value = 1.28
$ws.Range("K48").Value = '07/05/2023 17:12'
$ws.Range("L48").Value = 1.27
$ws.Range("M48").Value = '13/05/2023 16:37'
$ws.Range("N48").Value = 6.04
$ws.Range("O48").Value = '07/05/2023 17:12'
$ws.Range("P48").Value = 6.13
$ws.Range("Q48").Value = '13/05/2023 16:58'
$ws.Range("R48").Value = 10.34
$ws.Range("S48").Value = '07/05/2023 17:12'
$ws.Range("T48").Value = 11.56
$ws.Range("U48").Value = '13/05/2023 16:58'
$ws.Range("V48").Value = 'https://www.betexplorer.com/football/norway/eliteserien/molde-stromsgodset/d6G7EpKs/'
$ws.Range("F49").Value = 'Ham-Kam'
$ws.Range("G49").Value = 1
$ws.Range("H49").Value = 'Tromso'
$ws.Range("I49").Value = 2
$ws.Range("J49").Value = 2.32
$ws.Range("K49").Value = '07/05/2023 17:12'
$ws.Range("L49").Value = 2.01
$ws.Range("M49").Value = '13/05/2023 16:50'
$ws.Range("N49").Value = 3.28
$ws.Range("O49").Value = '07/05/2023 17:12'
$ws.Range("P49").Value = 3.36
$ws.Range("Q49").Value = '13/05/2023 16:50'
$ws.Range("R49").Value = 3.3
$ws.Range("S49").Value = '07/05/2023 17:12'
$ws.Range("T49").Value = 4.25
$ws.Range("U49").Value = '13/05/2023 16:57'
$ws.Range("V49").Value = 'https://www.betexplorer.com/football/norway/eliteserien/ham-kam-tromso/MDPuJQBQ/'
$ws.Range("F51").Value = 'Sandefjord'
$ws.Range("G51").Value = 1
$ws.Range("H51").Value = 'Viking'
$ws.Range("I51").Value = 2
$ws.Range("J51").Value = 3.58
$ws.Range("K51").Value = '08/05/2023 19:12'
$ws.Range("L51").Value = 3.8
$ws.Range("M51").Value = '13/05/2023 16:58'
$ws.Range("N51").Value = 3.74
$ws.Range("O51").Value = '08/05/2023 19:12'
$ws.Range("P51").Value = 3.77
$ws.Range("Q51").Value = '13/05/2023 16:58'
$ws.Range("R51").Value = 2.03
$ws.Range("S51").Value = '08/05/2023 19:12'
$ws.Range("T51").Value = 2
$ws.Range("U51").Value = '13/05/2023 16:57'
$ws.Range("V51").Value = 'https://www.betexplorer.com/football/norway/eliteserien/sandefjord-viking/vV5CD4Zm/'
$ws.Range("F167").Value = 'Stabaek'
$ws.Range("G167").Value = 0
$ws.Range("H167").Value = 'Brann'
$ws.Range("I167").Value = 1
$ws.Range("J167").Value = 3.78
$ws.Range("K167").Value = '04/09/2023 16:12'
$ws.Range("L167").Value = 3.94
$ws.Range("M167").Value = '17/09/2023 16:59'
$ws.Range("N167").Value = 3.87
$ws.Range("O167").Value = '04/09/2023 16:12'
$ws.Range("P167").Value = 3.9
$ws.Range("Q167").Value = '17/09/2023 16:51'
$ws.Range("R167").Value = 1.93
$ws.Range("S167").Value = '04/09/2023 16:12'
$ws.Range("T167").Value = 1.79
$ws.Range("U167").Value = '17/09/2023 16:59'
$ws.Range("V167").Value = 'https://www.betexplorer.com/football/norway/eliteserien/stabaek-brann/44DZhnBm/'
$ws.Range("F168").Value = 'Tromso'
$ws.Range("G168").Value = 2
$ws.Range("H168").Value = 'HamKam'
$ws.Range("I168").Value = 1
$ws.Range("J168").Value = 1.6
$ws.Range("K168").Value = '04/09/2023 16:12'
$ws.Range("L168").Value = 1.51
$ws.Range("M168").Value = '17/09/2023 16:51'
$ws.Range("N168").Value = 4.26
$ws.Range("O168").Value = '04/09/2023 16:12'
$ws.Range("P168").Value = 4.56
$ws.Range("Q168").Value = '17/09/2023 16:38'
$ws.Range("R168").Value = 5.48
$ws.Range("S168").Value = '04/09/2023 16:12'
$ws.Range("T168").Value = 6.89
$ws.Range("U168").Value = '17/09/2023 16:38'
$ws.Range("V168").Value = 'https://www.betexplorer.com/football/norway/eliteserien/tromso-ham-kam/bRCwh6Qg/'
$ws.Range("F169").Value = 'Rosenborg'
$ws.Range("G169").Value = 1
$ws.Range("H169").Value = 'Bodo/Glimt'
$ws.Range("I169").Value = 1
$ws.Range("J169").Value = 3.83
$ws.Range("K169").Value = '04/09/2023 16:12'
$ws.Range("L169").Value = 4.31
$ws.Range("M169").Value = '17/09/2023 16:59'
$ws.Range("N169").Value = 4.1
$ws.Range("O169").Value = '04/09/2023 16:12'
$ws.Range("P169").Value = 4.48
$ws.Range("Q169").Value = '17/09/2023 16:57'
$ws.Range("R169").Value = 1.86
$ws.Range("S169").Value = '04/09/2023 16:12'
$ws.Range("T169").Value = 1.65
$ws.Range("U169").Value = '17/09/2023 16:54'
$ws.Range("V169").Value = 'https://www.betexplorer.com/football/norway/eliteserien/rosenborg-bodo-glimt/6TXbalYJ/'
$ws.Range("F170").Value = 'Sarpsborg 08'
$ws.Range("G170").Value = 3
$ws.Range("H170").Value = 'Lillestrom'
$ws.Range("I170").Value = 1
$ws.Range("J170").Value = 2.05
$ws.Range("K170").Value = '04/09/2023 16:12'
$ws.Range("L170").Value = 1.94
$ws.Range("M170").Value = '17/09/2023 16:45'
$ws.Range("N170").Value = 4
$ws.Range("O170").Value = '04/09/2023 16:12'
$ws.Range("P170").Value = 4.2
$ws.Range("Q170").Value = '17/09/2023 16:57'
$ws.Range("R170").Value = 3.32
$ws.Range("S170").Value = '04/09/2023 16:12'
$ws.Range("T170").Value = 3.63
$ws.Range("U170").Value = '17/09/2023 16:57'
$ws.Range("V170").Value = 'https://www.betexplorer.com/football/norway/eliteserien/sarpsborg-08-lillestrom/QDEVgSes/'
$ws.Range("F173").Value = 'Viking'
$ws.Range("G173").Value = 4
$ws.Range("H173").Value = 'Sandefjord'
$ws.Range("I173").Value = 3
$ws.Range("J173").Value = 1.28
$ws.Range("K173").Value = '16/09/2023 17:13'
$ws.Range("L173").Value = 1.28
$ws.Range("M173").Value = '24/09/2023 16:40'
$ws.Range("N173").Value = 6.7
$ws.Range("O173").Value = '16/09/2023 17:13'
$ws.Range("P173").Value = 6.68
$ws.Range("Q173").Value = '24/09/2023 16:55'
$ws.Range("R173").Value = 9.140000000000001
$ws.Range("S173").Value = '16/09/2023 17:13'
$ws.Range("T173").Value = 9.109999999999999
$ws.Range("U173").Value = '24/09/2023 16:55'
$ws.Range("V173").Value = 'https://www.betexplorer.com/football/norway/eliteserien/viking-sandefjord/Y7JEkYLh/'
$ws.Range("F174").Value = 'Odd'
$ws.Range("G174").Value = 1
$ws.Range("H174").Value = 'Haugesund'
$ws.Range("I174").Value = 1
$ws.Range("J174").Value = 1.92
$ws.Range("K174").Value = '16/09/2023 17:13'
$ws.Range("L174").Value = 2.02
$ws.Range("M174").Value = '24/09/2023 16:26'
$ws.Range("N174").Value = 3.76
$ws.Range("O174").Value = '16/09/2023 17:13'
$ws.Range("P174").Value = 3.63
$ws.Range("Q174").Value = '24/09/2023 16:31'
$ws.Range("R174").Value = 3.96
$ws.Range("S174").Value = '16/09/2023 17:13'
$ws.Range("T174").Value = 3.88
$ws.Range("U174").Value = '24/09/2023 16:44'
$ws.Range("V174").Value = 'https://www.betexplorer.com/football/norway/eliteserien/odds-bk-haugesund/0fL6ifit/'
$ws.Range("F175").Value = 'Stromsgodset'
$ws.Range("G175").Value = 1
$ws.Range("H175").Value = 'Molde'
$ws.Range("I175").Value = 1
$ws.Range("J175").Value = 4.42
$ws.Range("K175").Value = '16/09/2023 17:13'
$ws.Range("L175").Value = 4.8
$ws.Range("M175").Value = '24/09/2023 16:34'
$ws.Range("N175").Value = 4.21
$ws.Range("O175").Value = '16/09/2023 17:13'
$ws.Range("P175").Value = 4.45
$ws.Range("Q175").Value = '24/09/2023 16:34'
$ws.Range("R175").Value = 1.73
$ws.Range("S175").Value = '16/09/2023 17:13'
$ws.Range("T175").Value = 1.64
$ws.Range("U175").Value = '24/09/2023 16:58'
$ws.Range("V175").Value = 'https://www.betexplorer.com/football/norway/eliteserien/stromsgodset-molde/fHKAjE6n/'
$ws.Range("F176").Value = 'HamKam'
$ws.Range("G176").Value = 1
$ws.Range("H176").Value = 'Sarpsborg 08'
$ws.Range("I176").Value = 1
$ws.Range("J176").Value = 3.11
$ws.Range("K176").Value = '17/09/2023 16:13'
$ws.Range("L176").Value = 3.62
$ws.Range("M176").Value = '24/09/2023 16:44'
$ws.Range("N176").Value = 3.93
$ws.Range("O176").Value = '17/09/2023 16:13'
$ws.Range("P176").Value = 4.41
$ws.Range("Q176").Value = '24/09/2023 16:49'
$ws.Range("R176").Value = 2.16
$ws.Range("S176").Value = '17/09/2023 16:13'
$ws.Range("T176").Value = 1.9
$ws.Range("U176").Value = '24/09/2023 16:44'
$ws.Range("V176").Value = 'https://www.betexplorer.com/football/norway/eliteserien/ham-kam-sarpsborg-08/K03bmruP/'
$ws.Range("F177").Value = 'Bodo/Glimt'
$ws.Range("G177").Value = 4
$ws.Range("H177").Value = 'Valerenga'
$ws.Range("I177").Value = 2
$ws.Range("J177").Value = 1.34
$ws.Range("K177").Value = '17/09/2023 18:43'
$ws.Range("L177").Value = 1.26
$ws.Range("M177").Value = '24/09/2023 16:56'
$ws.Range("N177").Value = 6.05
$ws.Range("O177").Value = '17/09/2023 18:43'
$ws.Range("P177").Value = 6.76
$ws.Range("Q177").Value = '24/09/2023 16:58'
$ws.Range("R177").Value = 7.63
$ws.Range("S177").Value = '17/09/2023 18:43'
$ws.Range("T177").Value = 10.37
$ws.Range("U177").Value = '24/09/2023 16:58'
$ws.Range("V177").Value = 'https://www.betexplorer.com/football/norway/eliteserien/bodo-glimt-valerenga/Yq5jk4AC/'
$ws.Range("F178").Value = 'Lillestrom'
$ws.Range("G178").Value = 3
$ws.Range("H178").Value = 'Rosenborg'
$ws.Range("I178").Value = 0
$ws.Range("J178").Value = 1.85
$ws.Range("K178").Value = '21/09/2023 10:43'
$ws.Range("L178").Value = 1.84
$ws.Range("M178").Value = '24/09/2023 16:59'
$ws.Range("N178").Value = 4.03
$ws.Range("O178").Value = '21/09/2023 10:43'
$ws.Range("P178").Value = 4.04
$ws.Range("Q178").Value = '24/09/2023 16:59'
$ws.Range("R178").Value = 3.95
$ws.Range("S178").Value = '21/09/2023 10:43'
$ws.Range("T178").Value = 4.21
$ws.Range("U178").Value = '24/09/2023 16:59'
$ws.Range("V178").Value = 'https://www.betexplorer.com/football/norway/eliteserien/lillestrom-rosenborg/OYTvdzqP/'
$ws.Range("F182").Value = 'Molde'
$ws.Range("G182").Value = 4
$ws.Range("H182").Value = 'Viking'
$ws.Range("I182").Value = 0
$ws.Range("J182").Value = 1.83
$ws.Range("K182").Value = '24/09/2023 16:13'
$ws.Range("L182").Value = 1.83
$ws.Range("M182").Value = '01/10/2023 16:54'
$ws.Range("N182").Value = 4.32
$ws.Range("O182").Value = '24/09/2023 16:13'
$ws.Range("P182").Value = 4.38
$ws.Range("Q182").Value = '01/10/2023 16:57'
$ws.Range("R182").Value = 3.82
$ws.Range("S182").Value = '24/09/2023 16:13'
$ws.Range("T182").Value = 3.91
$ws.Range("U182").Value = '01/10/2023 16:57'
$ws.Range("V182").Value = 'https://www.betexplorer.com/football/norway/eliteserien/molde-viking/QwbMSkcr/'
$ws.Range("F183").Value = 'Sandefjord'
$ws.Range("G183").Value = 1
$ws.Range("H183").Value = 'Valerenga'
$ws.Range("I183").Value = 2
$ws.Range("J183").Value = 3.17
$ws.Range("K183").Value = '24/09/2023 16:13'
$ws.Range("L183").Value = 2.95
$ws.Range("M183").Value = '01/10/2023 16:39'
$ws.Range("N183").Value = 3.62
$ws.Range("O183").Value = '24/09/2023 16:13'
$ws.Range("P183").Value = 3.56
$ws.Range("Q183").Value = '01/10/2023 16:52'
$ws.Range("R183").Value = 2.24
$ws.Range("S183").Value = '24/09/2023 16:13'
$ws.Range("T183").Value = 2.45
$ws.Range("U183").Value = '01/10/2023 15:46'
$ws.Range("V183").Value = 'https://www.betexplorer.com/football/norway/eliteserien/sandefjord-valerenga/vBPdYHYR/'
$ws.Range("F187").Value = 'Haugesund'
$ws.Range("G187").Value = 1
$ws.Range("H187").Value = 'Stromsgodset'
$ws.Range("I187").Value = 0
$ws.Range("J187").Value = 2.31
$ws.Range("K187").Value = '01/10/2023 18:42'
$ws.Range("L187").Value = 2.51
$ws.Range("M187").Value = '08/10/2023 16:58'
$ws.Range("N187").Value = 3.6
$ws.Range("O187").Value = '01/10/2023 18:42'
$ws.Range("P187").Value = 3.42
$ws.Range("Q187").Value = '08/10/2023 16:59'
$ws.Range("R187").Value = 3.04
$ws.Range("S187").Value = '01/10/2023 18:42'
$ws.Range("T187").Value = 2.98
$ws.Range("U187").Value = '08/10/2023 16:58'
$ws.Range("V187").Value = 'https://www.betexplorer.com/football/norway/eliteserien/haugesund-stromsgodset/jqEJlhya/'
$ws.Range("F188").Value = 'Molde'
$ws.Range("G188").Value = 1
$ws.Range("H188").Value = 'Bodo/Glimt'
$ws.Range("I188").Value = 3
$ws.Range("J188").Value = 2.22
$ws.Range("K188").Value = '01/10/2023 18:42'
$ws.Range("L188").Value = 2.13
$ws.Range("M188").Value = '08/10/2023 16:52'
$ws.Range("N188").Value = 3.83
$ws.Range("O188").Value = '01/10/2023 18:42'
$ws.Range("P188").Value = 3.94
$ws.Range("Q188").Value = '08/10/2023 16:55'
$ws.Range("R188").Value = 3.11
$ws.Range("S188").Value = '01/10/2023 18:42'
$ws.Range("T188").Value = 3.28
$ws.Range("U188").Value = '08/10/2023 16:51'
$ws.Range("V188").Value = 'https://www.betexplorer.com/football/norway/eliteserien/molde-bodo-glimt/UZCNmCj5/'
$ws.Range("F189").Value = 'Sandefjord'
$ws.Range("G189").Value = 0
$ws.Range("H189").Value = 'HamKam'
$ws.Range("I189").Value = 1
$ws.Range("J189").Value = 2.15
$ws.Range("K189").Value = '04/10/2023 18:12'
$ws.Range("L189").Value = 2.03
$ws.Range("M189").Value = '08/10/2023 16:52'
$ws.Range("N189").Value = 3.75
$ws.Range("O189").Value = '04/10/2023 18:12'
$ws.Range("P189").Value = 3.9
$ws.Range("Q189").Value = '08/10/2023 16:52'
$ws.Range("R189").Value = 3.31
$ws.Range("S189").Value = '04/10/2023 18:12'
$ws.Range("T189").Value = 3.57
$ws.Range("U189").Value = '08/10/2023 16:52'
$ws.Range("V189").Value = 'https://www.betexplorer.com/football/norway/eliteserien/sandefjord-ham-kam/zXFVojLH/'
$ws.Range("F190").Value = 'Sarpsborg 08'
$ws.Range("G190").Value = 5
$ws.Range("H190").Value = 'Rosenborg'
$ws.Range("I190").Value = 2
$ws.Range("J190").Value = 1.74
$ws.Range("K190").Value = '30/09/2023 17:13'
$ws.Range("L190").Value = 1.6
$ws.Range("M190").Value = '08/10/2023 16:53'
$ws.Range("N190").Value = 4.36
$ws.Range("O190").Value = '30/09/2023 17:13'
$ws.Range("P190").Value = 4.76
$ws.Range("Q190").Value = '08/10/2023 16:59'
$ws.Range("R190").Value = 4.39
$ws.Range("S190").Value = '30/09/2023 17:13'
$ws.Range("T190").Value = 5.12
$ws.Range("U190").Value = '08/10/2023 16:58'
$ws.Range("V190").Value = 'https://www.betexplorer.com/football/norway/eliteserien/sarpsborg-08-rosenborg/GUTaUNkp/'
$ws.Range("F191").Value = 'Stabaek'
$ws.Range("G191").Value = 1
$ws.Range("H191").Value = 'Lillestrom'
$ws.Range("I191").Value = 0
$ws.Range("J191").Value = 3.21
$ws.Range("K191").Value = '04/10/2023 18:12'
$ws.Range("L191").Value = 2.26
$ws.Range("M191").Value = '08/10/2023 16:51'
$ws.Range("N191").Value = 3.53
$ws.Range("O191").Value = '04/10/2023 18:12'
$ws.Range("P191").Value = 3.56
$ws.Range("Q191").Value = '08/10/2023 16:51'
$ws.Range("R191").Value = 2.29
$ws.Range("S191").Value = '04/10/2023 18:12'
$ws.Range("T191").Value = 3.27
$ws.Range("U191").Value = '08/10/2023 16:51'
$ws.Range("V191").Value = 'https://www.betexplorer.com/football/norway/eliteserien/stabaek-lillestrom/fLS3Ts5j/'
$ws.Range("F192").Value = 'Tromso'
$ws.Range("G192").Value = 1
$ws.Range("H192").Value = 'Aalesund'
$ws.Range("I192").Value = 2
$ws.Range("J192").Value = 1.37
$ws.Range("K192").Value = '30/09/2023 17:13'
$ws.Range("L192").Value = 1.29
$ws.Range("M192").Value = '08/10/2023 16:52'
$ws.Range("N192").Value = 5.21
$ws.Range("O192").Value = '30/09/2023 17:13'
$ws.Range("P192").Value = 5.99
$ws.Range("Q192").Value = '08/10/2023 16:52'
$ws.Range("R192").Value = 7.99
$ws.Range("S192").Value = '30/09/2023 17:13'
$ws.Range("T192").Value = 10.82
$ws.Range("U192").Value = '08/10/2023 16:59'
$ws.Range("V192").Value = 'https://www.betexplorer.com/football/norway/eliteserien/tromso-aalesund/ryMohRsi/'
$ws.Range("F197").Value = 'Aalesund'
$ws.Range("G197").Value = 3
$ws.Range("H197").Value = 'Sarpsborg 08'
$ws.Range("I197").Value = 2
$ws.Range("J197").Value = 3.9
$ws.Range("K197").Value = '09/10/2023 14:42'
$ws.Range("L197").Value = 3.58
$ws.Range("M197").Value = '22/10/2023 16:53'
$ws.Range("N197").Value = 4.25
$ws.Range("O197").Value = '09/10/2023 14:42'
$ws.Range("P197").Value = 4.21
$ws.Range("Q197").Value = '22/10/2023 16:53'
$ws.Range("R197").Value = 1.82
$ws.Range("S197").Value = '09/10/2023 14:42'
$ws.Range("T197").Value = 1.95
$ws.Range("U197").Value = '22/10/2023 16:53'
$ws.Range("V197").Value = 'https://www.betexplorer.com/football/norway/eliteserien/aalesund-sarpsborg-08/4jAfj5C3/'
$ws.Range("F198").Value = 'HamKam'
$ws.Range("G198").Value = 0
$ws.Range("H198").Value = 'Haugesund'
$ws.Range("I198").Value = 3
$ws.Range("J198").Value = 2.05
$ws.Range("K198").Value = '09/10/2023 14:42'
$ws.Range("L198").Value = 2.28
$ws.Range("M198").Value = '22/10/2023 16:55'
$ws.Range("N198").Value = 3.74
$ws.Range("O198").Value = '09/10/2023 14:42'
$ws.Range("P198").Value = 3.52
$ws.Range("Q198").Value = '22/10/2023 16:55'
$ws.Range("R198").Value = 3.51
$ws.Range("S198").Value = '09/10/2023 14:42'
$ws.Range("T198").Value = 3.29
$ws.Range("U198").Value = '22/10/2023 16:55'
$ws.Range("V198").Value = 'https://www.betexplorer.com/football/norway/eliteserien/ham-kam-haugesund/j7C6m3dM/'
$ws.Range("F199").Value = 'Rosenborg'
$ws.Range("G199").Value = 1
$ws.Range("H199").Value = 'Stabaek'
$ws.Range("I199").Value = 1
$ws.Range("J199").Value = 1.79
$ws.Range("K199").Value = '09/10/2023 14:42'
$ws.Range("L199").Value = 1.98
$ws.Range("M199").Value = '22/10/2023 16:56'
$ws.Range("N199").Value = 4.01
$ws.Range("O199").Value = '09/10/2023 14:42'
$ws.Range("P199").Value = 3.88
$ws.Range("Q199").Value = '22/10/2023 16:53'
$ws.Range("R199").Value = 4.27
$ws.Range("S199").Value = '09/10/2023 14:42'
$ws.Range("T199").Value = 3.75
$ws.Range("U199").Value = '22/10/2023 16:56'
$ws.Range("V199").Value = 'https://www.betexplorer.com/football/norway/eliteserien/rosenborg-stabaek/KxVUssKq/'
$ws.Range("F200").Value = 'Stromsgodset'
$ws.Range("G200").Value = 3
$ws.Range("H200").Value = 'Odd'
$ws.Range("I200").Value = 1
$ws.Range("J200").Value = 2.15
$ws.Range("K200").Value = '09/10/2023 14:42'
$ws.Range("L200").Value = 1.97
$ws.Range("M200").Value = '22/10/2023 16:53'
$ws.Range("N200").Value = 3.63
$ws.Range("O200").Value = '09/10/2023 14:42'
$ws.Range("P200").Value = 3.56
$ws.Range("Q200").Value = '22/10/2023 16:54'
$ws.Range("R200").Value = 3.42
$ws.Range("S200").Value = '09/10/2023 14:42'
$ws.Range("T200").Value = 4.14
$ws.Range("U200").Value = '22/10/2023 16:53'
$ws.Range("V200").Value = 'https://www.betexplorer.com/football/norway/eliteserien/stromsgodset-odds-bk/EaUYt1Zk/'

# ===========================================================================
# 2) Append newly scraped fixtures as rows 211-216, copying the row format
#    (styles/number formats) from the last existing row (210) first.
# ===========================================================================
$ws.Range("A210:V210").Copy()
$ws.Range("A211:V216").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- row 211 ---
$ws.Range("A211").Value = 210
$ws.Range("B211").Value = 'norway'
$ws.Range("C211").Value = 'eliteserien'
$ws.Range("D211").Formula = "'2023"
$ws.Range("E211").Value = 45235.70833333334
$ws.Range("F211").Value = 'Brann'
$ws.Range("G211").Value = 2
$ws.Range("H211").Value = 'Odd'
$ws.Range("I211").Value = 1
$ws.Range("J211").Value = 1.34
$ws.Range("K211").Value = '29/10/2023 17:13'
$ws.Range("L211").Value = 1.14
$ws.Range("M211").Value = '05/11/2023 16:36'
$ws.Range("N211").Value = 5.76
$ws.Range("O211").Value = '29/10/2023 17:13'
$ws.Range("P211").Value = 9.7
$ws.Range("Q211").Value = '05/11/2023 16:56'
$ws.Range("R211").Value = 8.05
$ws.Range("S211").Value = '29/10/2023 17:13'
$ws.Range("T211").Value = 18.39
$ws.Range("U211").Value = '05/11/2023 16:56'
$ws.Range("V211").Value = 'https://www.betexplorer.com/football/norway/eliteserien/brann-odds-bk/2LnLhbY1/'
# --- row 212 ---
$ws.Range("A212").Value = 211
$ws.Range("B212").Value = 'norway'
$ws.Range("C212").Value = 'eliteserien'
$ws.Range("D212").Formula = "'2023"
$ws.Range("E212").Value = 45235.70833333334
$ws.Range("F212").Value = 'Lillestrom'
$ws.Range("G212").Value = 1
$ws.Range("H212").Value = 'Haugesund'
$ws.Range("I212").Value = 0
$ws.Range("J212").Value = 1.53
$ws.Range("K212").Value = '30/10/2023 18:42'
$ws.Range("L212").Value = 1.87
$ws.Range("M212").Value = '05/11/2023 16:59'
$ws.Range("N212").Value = 4.64
$ws.Range("O212").Value = '30/10/2023 18:42'
$ws.Range("P212").Value = 3.81
$ws.Range("Q212").Value = '05/11/2023 16:59'
$ws.Range("R212").Value = 5.76
$ws.Range("S212").Value = '30/10/2023 18:42'
$ws.Range("T212").Value = 4.29
$ws.Range("U212").Value = '05/11/2023 16:59'
$ws.Range("V212").Value = 'https://www.betexplorer.com/football/norway/eliteserien/lillestrom-haugesund/URv7G1Q8/'
# --- row 213 ---
$ws.Range("A213").Value = 212
$ws.Range("B213").Value = 'norway'
$ws.Range("C213").Value = 'eliteserien'
$ws.Range("D213").Formula = "'2023"
$ws.Range("E213").Value = 45235.70833333334
$ws.Range("F213").Value = 'Rosenborg'
$ws.Range("G213").Value = 3
$ws.Range("H213").Value = 'Molde'
$ws.Range("I213").Value = 1
$ws.Range("J213").Value = 3.55
$ws.Range("K213").Value = '30/10/2023 21:13'
$ws.Range("L213").Value = 3.65
$ws.Range("M213").Value = '05/11/2023 16:59'
$ws.Range("N213").Value = 3.87
$ws.Range("O213").Value = '30/10/2023 21:13'
$ws.Range("P213").Value = 4.06
$ws.Range("Q213").Value = '05/11/2023 16:59'
$ws.Range("R213").Value = 2
$ws.Range("S213").Value = '30/10/2023 21:13'
$ws.Range("T213").Value = 1.97
$ws.Range("U213").Value = '05/11/2023 16:59'
$ws.Range("V213").Value = 'https://www.betexplorer.com/football/norway/eliteserien/rosenborg-molde/CjUBFLuF/'
# --- row 214 ---
$ws.Range("A214").Value = 213
$ws.Range("B214").Value = 'norway'
$ws.Range("C214").Value = 'eliteserien'
$ws.Range("D214").Formula = "'2023"
$ws.Range("E214").Value = 45235.70833333334
$ws.Range("F214").Value = 'Sarpsborg 08'
$ws.Range("G214").Value = 3
$ws.Range("H214").Value = 'Valerenga'
$ws.Range("I214").Value = 2
$ws.Range("J214").Value = 1.75
$ws.Range("K214").Value = '29/10/2023 19:43'
$ws.Range("L214").Value = 1.91
$ws.Range("M214").Value = '05/11/2023 16:51'
$ws.Range("N214").Value = 4.33
$ws.Range("O214").Value = '29/10/2023 19:43'
$ws.Range("P214").Value = 4.22
$ws.Range("Q214").Value = '05/11/2023 16:51'
$ws.Range("R214").Value = 4.34
$ws.Range("S214").Value = '29/10/2023 19:43'
$ws.Range("T214").Value = 3.72
$ws.Range("U214").Value = '05/11/2023 16:51'
$ws.Range("V214").Value = 'https://www.betexplorer.com/football/norway/eliteserien/sarpsborg-08-valerenga/jsTFEueL/'
# --- row 215 ---
$ws.Range("A215").Value = 214
$ws.Range("B215").Value = 'norway'
$ws.Range("C215").Value = 'eliteserien'
$ws.Range("D215").Formula = "'2023"
$ws.Range("E215").Value = 45235.70833333334
$ws.Range("F215").Value = 'Stabaek'
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = 'Bodo/Glimt'
$ws.Range("I215").Value = 4
$ws.Range("J215").Value = 4.76
$ws.Range("K215").Value = '01/11/2023 16:42'
$ws.Range("L215").Value = 5.73
$ws.Range("M215").Value = '05/11/2023 16:59'
$ws.Range("N215").Value = 4.57
$ws.Range("O215").Value = '01/11/2023 16:42'
$ws.Range("P215").Value = 4.97
$ws.Range("Q215").Value = '05/11/2023 16:59'
$ws.Range("R215").Value = 1.63
$ws.Range("S215").Value = '01/11/2023 16:42'
$ws.Range("T215").Value = 1.52
$ws.Range("U215").Value = '05/11/2023 16:56'
$ws.Range("V215").Value = 'https://www.betexplorer.com/football/norway/eliteserien/stabaek-bodo-glimt/b7SJDaAR/'
# --- row 216 ---
$ws.Range("A216").Value = 215
$ws.Range("B216").Value = 'norway'
$ws.Range("C216").Value = 'eliteserien'
$ws.Range("D216").Formula = "'2023"
$ws.Range("E216").Value = 45235.80208333334
$ws.Range("F216").Value = 'Tromso'
$ws.Range("G216").Value = 0
$ws.Range("H216").Value = 'Stromsgodset'
$ws.Range("I216").Value = 1
$ws.Range("J216").Value = 1.56
$ws.Range("K216").Value = '30/10/2023 21:13'
$ws.Range("L216").Value = 1.46
$ws.Range("M216").Value = '05/11/2023 18:14'
$ws.Range("N216").Value = 4.45
$ws.Range("O216").Value = '30/10/2023 21:13'
$ws.Range("P216").Value = 4.6
$ws.Range("Q216").Value = '05/11/2023 19:13'
$ws.Range("R216").Value = 5.83
$ws.Range("S216").Value = '30/10/2023 21:13'
$ws.Range("T216").Value = 7.55
$ws.Range("U216").Value = '05/11/2023 19:13'
$ws.Range("V216").Value = 'https://www.betexplorer.com/football/norway/eliteserien/tromso-stromsgodset/0GEo9JIr/'

# Re-apply the row format after writing the values: forcing column D's
# "2023" season value to be text (via a leading apostrophe) makes Excel
# tag the cell with a transient quote-prefix style, so re-pasting the
# reference formatting cleans that back up to match the rest of the sheet.
$ws.Range("A210:V210").Copy()
$ws.Range("A211:V216").PasteSpecial(-4122)
$excel.CutCopyMode = 0
